$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 stays GLD / StreetTRACKS Gold Shares - update metrics
$ws.Range("D2").Value = 387.13
$ws.Range("E2").Value = 56.3
$ws.Range("F2").Value = 1.05
$ws.Range("I2").Value = 83
$ws.Range("K2").Value = 67.7
$ws.Range("M2").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N2").Value = 54.85170003294819
$ws.Range("O2").Value = "⚪ 중립 구간"

# Row 3 now becomes Gold Feb 26 / GC=F (swapped with old row 4 contents)
$ws.Range("B3").Value = "Gold Feb 26"
$ws.Range("C3").Value = "GC=F"
$ws.Range("D3").Value = 4229.4
$ws.Range("E3").Value = 54.8
$ws.Range("F3").Value = 1.54
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 63
$ws.Range("I3").Value = 80
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 66.5
$ws.Range("M3").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N3").Value = 54.85170003294819
$ws.Range("O3").Value = "⚪ 중립 구간"

# Row 4 now becomes Newmont Corporation / NEM (swapped with old row 3 contents)
$ws.Range("B4").Value = "Newmont Corporation"
$ws.Range("C4").Value = "NEM"
$ws.Range("D4").Value = 90.72
$ws.Range("E4").Value = 52.7
$ws.Range("F4").Value = 0.22
$ws.Range("G4").Value = 60
$ws.Range("H4").Value = 80
$ws.Range("I4").Value = 80
$ws.Range("J4").Value = 86
$ws.Range("K4").Value = 66.5
$ws.Range("M4").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N4").Value = 54.85170003294819
$ws.Range("O4").Value = "⚪ 중립 구간"
